$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.255.18"
$ws.Range("E2").Value = "  +0.11%  "

$ws.Range("D3").Value = "1.596.64"
$ws.Range("E3").Value = "  +0.57%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.18%  "

$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  +0.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0853"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.76%  "

$ws.Range("D12").Value = "1.822.76"
$ws.Range("E12").Value = "  +0.71%  "

$ws.Range("D13").Value = "1.594.36"
$ws.Range("E13").Value = "  +0.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.42%  "

$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.503"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.48%  "

$ws.Range("D17").Value = "26.250.90"
$ws.Range("E17").Value = "  +0.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +7.51%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.19%  "

$ws.Range("E20").Value = "  -0.46%  "

$ws.Range("E22").Value = "  -0.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.43%  "

$ws.Range("E24").Value = "  +1.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.93%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.112"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.35%  "

$ws.Range("E30").Value = "  -0.52%  "

$ws.Range("E31").Value = "  -0.23%  "

$ws.Range("E32").Value = "  +0.66%  "

$ws.Range("D33").Value = "1.464.25"
$ws.Range("E33").Value = "  +3.00%  "

$ws.Range("E35").Value = "  -0.18%  "

$ws.Range("E36").Value = "  +0.60%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.567"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.43%  "

$ws.Range("E38").Value = "  -1.11%  "

$ws.Range("E39").Value = "  -0.25%  "

$ws.Range("E40").Value = "  -2.34%  "

$ws.Range("E42").Value = "  +2.29%  "

$ws.Range("E43").Value = "  -1.54%  "

$ws.Range("D44").Value = "1.734.47"
$ws.Range("E44").Value = "  +0.79%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.755"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.49%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.50"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.99%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.44%  "

$ws.Range("E48").Value = "  -0.84%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0500"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.02%  "

$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.04%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0946"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.33%  "
